$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Disciplinary")
Write-Output $ws.Name
